$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 83: fix up the timestamp in column A only (rest of the row is unchanged) ---
$ws.Cells.Item(83, 1).Value = 45456.2916666667

# --- Row 84: brand new data row appended by the R script ---

# Column A is a date/time serial; copy A83's format (custom "yyyy-mm-dd hh:mm:ss")
# onto A84 so it reuses the existing style instead of Excel minting a new one.
$ws.Cells.Item(83, 1).Copy() | Out-Null
$ws.Cells.Item(84, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(84, 1).Value = 45457.6415740741

$ws.Cells.Item(84, 2).Value = 20100
$ws.Cells.Item(84, 3).Value = 6.40000009536743
$ws.Cells.Item(84, 4).Value = 5.98000001907349
$ws.Cells.Item(84, 5).Value = 6.17999982833862
$ws.Cells.Item(84, 6).Value = 6

# Column G (adj_close) is stored as text "6" in the source data, even though it
# looks numeric -- force text so it lands in sharedStrings rather than as a number.
$g84 = $ws.Cells.Item(84, 7)
$g84.NumberFormat = "@"
$g84.Value = "6"
$g84.Style = "Normal"

$ws.Cells.Item(84, 8).Value = "PAL.MI"
